$d = $word.ActiveDocument

# 1. Update the description of the first row ("01") in the issues table:
#    "Os dois CRUD's da release 02 não foram totalmente implementados."
#    becomes
#    "O CRUD da release 01 e 02 não foram totalmente implementados, falta permitir alterações."
$d.Content.Find.Execute(
    "Os dois CRUD" + [char]0x2019 + "s da release 02 não foram totalmente implementados.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "O CRUD da release 01 e 02 não foram totalmente implementados, falta permitir alterações.",
    2) | Out-Null

# 2. Clear out the "02" row (ID + Descrição) in the same table — the row that
#    described "Os testes da release 02 não foram implementados." is emptied
#    (row/cells are kept, but the text + its run formatting is removed).
$table = $d.Tables.Item(3)

$idCell = $table.Cell(3, 1)
$idRange = $d.Range($idCell.Range.Start, $idCell.Range.End)
$idRange.Find.Execute("02", $false, $false, $false, $false, $false, $true, 0, $false, "", 1) | Out-Null

$descCell = $table.Cell(3, 2)
$descRange = $d.Range($descCell.Range.Start, $descCell.Range.End)
$descRange.Find.Execute("Os testes da release 02 não foram implementados.", $false, $false, $false, $false, $false, $true, 0, $false, "", 1) | Out-Null

# 3. Extend the cronograma paragraph with the new trailing sentence about tests.
$d.Content.Find.Execute(
    "O cronograma foi refatorado e está mais próximo da realidade do projeto, além de implementar os recursos do projeto e acompanhamento da conclusão das tarefas. Ademais, agora é possível criar, consultar e apagar kitnets e repúblicas.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "O cronograma foi refatorado e está mais próximo da realidade do projeto, além de implementar os recursos do projeto e acompanhamento da conclusão das tarefas. Ademais, agora é possível criar, consultar e apagar kitnets e repúblicas, bem como seus testes foram implementados.",
    2) | Out-Null
